$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New iteration row (8th iteration) appended to the existing A:C table ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "33/200"
$ws.Range("C9").Value = 2200

# --- New "split" identifier table (breed/eng) in columns E:H ---
$ws.Range("E1").Value = "Iteration "
$ws.Range("F1").Value = "difference"
$ws.Range("G1").Value = "TP"
$ws.Range("H1").Value = "training data "

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "69/200"
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 0

# --- Leave the selection where the author left off (ready for 9th iteration) ---
$ws.Range("H10").Select() | Out-Null
